# Add a new "abbreviations" worksheet between "similar_facnames" and
# "sending_agencies", populate it with the abbreviation -> correction
# lookup table, and leave the new tab as the active/selected sheet.

$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item(1)

# Insert the new sheet right after "similar_facnames" (i.e. before the
# current second sheet, "sending_agencies").
$ws = $wb.Worksheets.Add($null, $firstSheet)
$ws.Name = "abbreviations"

# The abbreviation / correction pairs (everything except the header row,
# which gets typed in afterwards).
$data = @(
    @("INST", "INSTITUTE"),
    @("CTR", "CENTER"),
    @("SCH", "SCHOOL"),
    @("SCH.", "SCHOOL"),
    @("INC.", "INC"),
    @("LLC.", "LLC"),
    @("DCC", "DAY CARE CENTER"),
    @("DAYCARE", "DAY CARE CENTER"),
    @("NONMEDICAID", "NON-MEDICAID"),
    @("NON MEDICAID", "NON-MEDICAID"),
    @("ECLC", "EARLY CHILDHOOD CENTER?"),
    @("CORP.", "CORP"),
    @("KID’S", "KIDS"),
    @("CCC", "CHILD CARE CENTER OR DAY CARE CENTER"),
    @("SRVCE", "SERVICE"),
    @("SER", "SERVICE"),
    @("SVS", "SERVICES"),
    @("&", "AND"),
    @("ST", "STREET"),
    @("ADMIN", "ADMINISTRATION"),
    @("ASSOC", "ASSOCIATION"),
    @("E OR E.", "EAST"),
    @("W OR W.", "WEST")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Header row gets added last (pushing the data down a row) -- this keeps
# the shared-string order matching how the sheet was actually built.
$ws.Rows.Item(1).Insert()
$ws.Cells.Item(1, 1).Value = "Seen on table"
$ws.Cells.Item(1, 2).Value = "Possible Correction"

$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# Restore the selection on the original first sheet, then leave the new
# "abbreviations" tab selected/active with B1 highlighted.
$firstSheet.Range("B11").Select() | Out-Null
$ws.Range("B1").Select() | Out-Null
$ws.Activate() | Out-Null
